$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column "Lead" before the existing last column (Description) ---
$ws.Columns.Item(9).Insert()
$ws.Columns.Item(9).ColumnWidth = 9.67

# Header
$ws.Range("I1").Value = "Lead"

# Data rows
$ws.Range("I2").Value = "Regina"
$ws.Range("I3").Value = "LeadMandatory"

# --- Formatting pasted-in look for the new Lead values (Arial, dark gray) ---
$ws.Range("I2:I3").ClearFormats()
$ws.Range("I2:I3").Font.Name = "Arial"
$ws.Range("I2:I3").Font.Color = 2696481

# Row 3's Lead value wraps onto two lines
$ws.Range("I3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 28.5

# --- Selection moves to the new Lead cell for row 2 ---
$ws.Range("I2").Select()
